$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Förändrad" (Changed) date column C for rows 2-5 from serial 45224 to 45233
$ws.Range("C2").Value = 45233
$ws.Range("C3").Value = 45233
$ws.Range("C4").Value = 45233
$ws.Range("C5").Value = 45233
